# Scheduled market-price refresh: updates the cached currentAveragePrice /
# LevePrice / LeveProfit columns (H:N) on affected leve rows across all
# eight job sheets, per the latest Kujata data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 408.27585
$ws.Range("I33").Value = 371.56522
$ws.Range("J33").Value = 549
$ws.Range("K33").Value = 371.56522
$ws.Range("L33").Value = 549
$ws.Range("M33").Value = -142.56522
$ws.Range("N33").Value = -1007
$ws.Range("H58").Value = 1031.8572
$ws.Range("J58").Value = 2080.875
$ws.Range("L58").Value = 6242.625
$ws.Range("N58").Value = -6542.625
$ws.Range("H76").Value = 4759.5
$ws.Range("I76").Value = 5778
$ws.Range("K76").Value = 5778
$ws.Range("M76").Value = -5463
$ws.Range("H79").Value = 4759.5
$ws.Range("I79").Value = 5778
$ws.Range("K79").Value = 5778
$ws.Range("M79").Value = -4686
$ws.Range("H129").Value = 815.55554
$ws.Range("I129").Value = 575.2857
$ws.Range("J129").Value = 873.5517
$ws.Range("K129").Value = 1725.8571
$ws.Range("L129").Value = 2620.6551
$ws.Range("M129").Value = 3274.1429
$ws.Range("N129").Value = -12620.6551

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3480.2988
$ws.Range("I32").Value = 3350.4353
$ws.Range("K32").Value = 3350.4353
$ws.Range("M32").Value = -3063.4353
$ws.Range("H61").Value = 50001136
$ws.Range("I61").Value = 62500948
$ws.Range("K61").Value = 62500948
$ws.Range("M61").Value = -62500736
$ws.Range("H63").Value = 2441.4736
$ws.Range("I63").Value = 2158.8
$ws.Range("J63").Value = 2755.5557
$ws.Range("K63").Value = 2158.8
$ws.Range("L63").Value = 2755.5557
$ws.Range("M63").Value = -1472.8
$ws.Range("N63").Value = -4127.5557
$ws.Range("H66").Value = 2441.4736
$ws.Range("I66").Value = 2158.8
$ws.Range("J66").Value = 2755.5557
$ws.Range("K66").Value = 10794
$ws.Range("L66").Value = 13777.7785
$ws.Range("M66").Value = -7362
$ws.Range("N66").Value = -20641.7785
$ws.Range("H122").Value = 1472.7838
$ws.Range("I122").Value = 1252.931
$ws.Range("K122").Value = 3758.793
$ws.Range("M122").Value = -1308.793
$ws.Range("H130").Value = 19952.666
$ws.Range("J130").Value = 19952.666
$ws.Range("L130").Value = 19952.666
$ws.Range("N130").Value = -29992.666
$ws.Range("H136").Value = 50001136
$ws.Range("I136").Value = 62500948
$ws.Range("K136").Value = 187502844
$ws.Range("M136").Value = -187500294
$ws.Range("H139").Value = 28641.666
$ws.Range("J139").Value = 28641.666
$ws.Range("L139").Value = 28641.666
$ws.Range("N139").Value = -38921.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H87").Value = 65000
$ws.Range("J87").Value = 65000
$ws.Range("L87").Value = 65000
$ws.Range("N87").Value = -67496
$ws.Range("H90").Value = 65000
$ws.Range("J90").Value = 65000
$ws.Range("L90").Value = 195000
$ws.Range("N90").Value = -207480
$ws.Range("H105").Value = 83334650
$ws.Range("I105").Value = 90910350
$ws.Range("K105").Value = 90910350
$ws.Range("M105").Value = -90908603

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1658.4762
$ws.Range("I58").Value = 1404
$ws.Range("J58").Value = 2740
$ws.Range("K58").Value = 1404
$ws.Range("L58").Value = 2740
$ws.Range("M58").Value = -1201
$ws.Range("N58").Value = -3146
$ws.Range("H62").Value = 13335513
$ws.Range("I62").Value = 2280
$ws.Range("K62").Value = 2280
$ws.Range("M62").Value = -1656
$ws.Range("H65").Value = 13335513
$ws.Range("I65").Value = 2280
$ws.Range("K65").Value = 11400
$ws.Range("M65").Value = -8280
$ws.Range("H86").Value = 4478872.5
$ws.Range("I86").Value = 11142027
$ws.Range("J86").Value = 36769.332
$ws.Range("K86").Value = 11142027
$ws.Range("L86").Value = 36769.332
$ws.Range("M86").Value = -11140904
$ws.Range("N86").Value = -39015.332
$ws.Range("H89").Value = 4478872.5
$ws.Range("I89").Value = 11142027
$ws.Range("J89").Value = 36769.332
$ws.Range("K89").Value = 55710135
$ws.Range("L89").Value = 183846.66
$ws.Range("M89").Value = -55704519
$ws.Range("N89").Value = -195078.66
$ws.Range("H108").Value = 28275.334
$ws.Range("J108").Value = 28275.334
$ws.Range("L108").Value = 28275.334
$ws.Range("N108").Value = -35955.334
$ws.Range("H122").Value = 4707.154
$ws.Range("I122").Value = 5176
$ws.Range("J122").Value = 1112.6666
$ws.Range("K122").Value = 15528
$ws.Range("L122").Value = 3337.9998
$ws.Range("M122").Value = -13078
$ws.Range("N122").Value = -8237.9998
$ws.Range("H130").Value = 32966.668
$ws.Range("J130").Value = 32966.668
$ws.Range("L130").Value = 32966.668
$ws.Range("N130").Value = -43006.668
$ws.Range("H136").Value = 1658.4762
$ws.Range("I136").Value = 1404
$ws.Range("J136").Value = 2740
$ws.Range("K136").Value = 4212
$ws.Range("L136").Value = 8220
$ws.Range("M136").Value = -1662
$ws.Range("N136").Value = -13320

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 2590.8
$ws.Range("J55").Value = 3125
$ws.Range("L55").Value = 9375
$ws.Range("N55").Value = -9729
$ws.Range("H131").Value = 14926565
$ws.Range("J131").Value = 1290.5167
$ws.Range("L131").Value = 3871.550099999999
$ws.Range("N131").Value = -13951.5501
$ws.Range("H132").Value = 950.6316
$ws.Range("J132").Value = 1139.8
$ws.Range("L132").Value = 10258.2
$ws.Range("N132").Value = -15318.2
$ws.Range("H140").Value = 3072.375
$ws.Range("I140").Value = 2085.0908
$ws.Range("K140").Value = 6255.2724
$ws.Range("M140").Value = -1075.2724

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 112502350
$ws.Range("I70").Value = 83336470
$ws.Range("J70").Value = 200000000
$ws.Range("K70").Value = 83336470
$ws.Range("L70").Value = 200000000
$ws.Range("M70").Value = -83336200
$ws.Range("N70").Value = -200000540
$ws.Range("H73").Value = 112502350
$ws.Range("I73").Value = 83336470
$ws.Range("J73").Value = 200000000
$ws.Range("K73").Value = 83336470
$ws.Range("L73").Value = 200000000
$ws.Range("M73").Value = -83335534
$ws.Range("N73").Value = -200001872

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H128").Value = 75000
$ws.Range("J128").Value = 75000
$ws.Range("L128").Value = 75000
$ws.Range("N128").Value = -84960

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 62501628
$ws.Range("I126").Value = 71429576
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 214288728
$ws.Range("L126").Value = 18000
$ws.Range("M126").Value = -214286258
$ws.Range("N126").Value = -22940
$ws.Range("H136").Value = 885.8570999999999
$ws.Range("I136").Value = 773.7059
$ws.Range("K136").Value = 2321.1177
$ws.Range("M136").Value = 228.8822999999998
